# Update the "Förändrad" (Changed) date column (C) for rows 2-15.
# The stored date serial number moves from 45180 (2023-09-11) to
# 45181 (2023-09-12) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45181
}
